$d = $word.ActiveDocument

# The two "Title" date/time paragraphs, the fuel-level label, the inline
# picture and its filename caption are collapsed into a single Title
# paragraph carrying the new timestamp.
$d.Content.Find.Execute("2024/08/20 Hora 09:47", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024/09/16 Hora 13:13", 2)

# Drop every other paragraph (second date title, fuel-level text, the
# picture paragraph, and the filename caption), walking back-to-front so
# indices stay valid as paragraphs disappear.
$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 2; $i--) {
    $paras.Item($i).Range.Delete()
}
